# Apply the changes described by the diff:
# 1. Add a new worksheet "Thông tin GV" after the existing sheet.
# 2. Move the "GV info" columns (G:K) from sheet1 into the new sheet,
#    inserting a new "Ngày bắt đầu" column before "Ngày BCCK".
# 3. Clear the now-empty G:K block on sheet1.
# 4. Update selections on both sheets and make the new sheet the active tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Create the new worksheet right after the first one ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Thông tin GV"

# --- Copy the teacher-info block (values + styles) to the new sheet ---
$ws1.Range("G1:K2").Copy($ws2.Range("A1"))

# --- Insert a new column at E (shifts old E:E -> F:F), matching column D/E layout ---
$ws2.Columns.Item(5).Insert(-4161)

# Give the new column E the same style as the (now shifted) F column, then set its values
$ws2.Range("F1:F2").Copy($ws2.Range("E1:E2"))
$ws2.Range("E1").Value = "Ngày bắt đầu"

# --- Column widths for the new sheet ---
$ws2.Columns.Item(1).ColumnWidth = 15
$ws2.Columns.Item(2).ColumnWidth = 10.714285714285714
$ws2.Columns.Item(3).ColumnWidth = 11.428571428571429
$ws2.Columns.Item(4).ColumnWidth = 14
$ws2.Columns.Item(5).ColumnWidth = 14
$ws2.Columns.Item(6).ColumnWidth = 16.142857142857142

# --- Remove the now-duplicated teacher-info block from sheet1 ---
$ws1.Range("G1:K2").Clear()

# --- Selections ---
$ws1.Range("G1:K2").Select()
$ws2.Range("F8").Select()

# --- Make the new sheet the active tab ---
$ws2.Activate()
